# The "Förändrad" (Changed) date in column C for every data row (2-344)
# is bumped by one day, from 2023-09-12 (serial 45181) to 2023-09-13
# (serial 45182). Update the whole column range in one shot, preserving
# the existing date number format / style already applied to those cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C344").Value = 45182
